$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ArrayForm")

# --- "Expanding Array Formulas" index/label ---
$ws.Range("A20").Value = "Expanding Array Formulas:"

# --- First expanding example: A21:A24 + C21:C24 -> E21:F24 ---
$ws.Range("A21").Value = 1
$ws.Range("A22").Value = 2
$ws.Range("A23").Value = 3
$ws.Range("A24").Value = 4
$ws.Range("C21").Value = 5
$ws.Range("C22").Value = 6
$ws.Range("C23").Value = 7
$ws.Range("C24").Value = 8
$ws.Range("E21:F24").FormulaArray = "=A21:A24+C21:C24"

# --- Second expanding example: A21:A24 * A30:D30 -> F28:I31 ---
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 4
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = 6
$ws.Range("C30").Value = 7
$ws.Range("D30").Value = 8
$ws.Range("F28:I31").FormulaArray = "=A21:A24*A30:D30"

# --- Third expanding example: A28:D28 + A30:D30 -> A32:D33 ---
$ws.Range("A32:D33").FormulaArray = "=A28:D28+A30:D30"

# --- Selection / active sheet bookkeeping ---
$ws.Range("L24").Select() | Out-Null

$wsSref = $wb.Worksheets.Item("sref")
$wsSref.Activate() | Out-Null
$wsSref.Range("A8").Select() | Out-Null
